$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the
#    "        month += 12;" paragraph to the very start of the
#    document (right before the first run of the first paragraph).
#
#    Bookmarks.Add("_GoBack", <range>) re-defines/moves an existing
#    bookmark of that name. A genuinely collapsed Range whose Start
#    and End are both 0 cannot be passed directly (it gets treated as
#    "no range supplied" and falls back to selecting the whole first
#    paragraph), so instead we insert a one-character placeholder at
#    the very start, wrap the bookmark tightly around just that
#    character, then delete the character - which correctly collapses
#    the bookmark down to a zero-length bookmark at position 0.
# ------------------------------------------------------------------

$startRng = $d.Range(0, 0)
$startRng.InsertBefore("X")

$bmRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range(0, 1).Text = ""

# ------------------------------------------------------------------
# 2) Mark the "Normal (Web)" style as a Quick Style (adds <w:qFormat/>
#    to its style definition).
# ------------------------------------------------------------------

$webStyle = $d.Styles("Normal (Web)")
$webStyle.QuickStyle = $true
